# VorlageWochenbericht.xlsx edit:
#  - Fix typo "Funktionsberreich" -> "Funktionsbereich" in the instructional
#    text held in cell A1 (a shared string).
#  - Leave the last worksheet selection on B12 (matches the saved cursor
#    position recorded in the workbook when it was re-saved).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$a1 = $ws.Range("A1")
$text = $a1.Value2
$fixed = $text -replace "Funktionsberreich", "Funktionsbereich"
$a1.Value = $fixed

$ws.Range("B12").Select()
